# Insert two new data rows into the Puerro price sheet.
#
# The sheet holds one observation per row (rows 2..200, header on row 1).
# The edit inserts a new row at position 86 (shifting the former rows
# 86..200 down to 87..201) and a second new row at position 124 of the
# *post-first-insert* numbering (shifting former rows 124..201 down to
# 125..202). The sheet dimension grows from A1:R200 to A1:R202
# automatically once the rows are inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert row 86 and populate it -----------------------------------
$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value2  = 10
$ws.Cells.Item(86, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(86, 3).Value2  = 'La Araucanía'
$ws.Cells.Item(86, 4).Value2  = 44679
$ws.Cells.Item(86, 5).Value2  = 9
$ws.Cells.Item(86, 6).Value2  = 100112005
$ws.Cells.Item(86, 7).Value2  = 'Puerro'
$ws.Cells.Item(86, 8).Value2  = 'Azul de Maquehue'
$ws.Cells.Item(86, 9).Value2  = 'Primera'
$ws.Cells.Item(86, 10).Value2 = 40
$ws.Cells.Item(86, 11).Value2 = 12000
$ws.Cells.Item(86, 12).Value2 = 12000
$ws.Cells.Item(86, 13).Value2 = 12000
$ws.Cells.Item(86, 14).Value2 = '$/docena de paquetes'
$ws.Cells.Item(86, 15).Value2 = 'Provincia de Cautín'
$ws.Cells.Item(86, 16).Value2 = 1000
$ws.Cells.Item(86, 17).Value2 = 12
$ws.Cells.Item(86, 18).Value2 = 'Hortaliza'

# --- Insert row 124 (post first-insert numbering) and populate it ----
$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value2  = 10
$ws.Cells.Item(124, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(124, 3).Value2  = 'La Araucanía'
$ws.Cells.Item(124, 4).Value2  = 44680
$ws.Cells.Item(124, 5).Value2  = 9
$ws.Cells.Item(124, 6).Value2  = 100112005
$ws.Cells.Item(124, 7).Value2  = 'Puerro'
$ws.Cells.Item(124, 8).Value2  = 'Azul de Maquehue'
$ws.Cells.Item(124, 9).Value2  = 'Primera'
$ws.Cells.Item(124, 10).Value2 = 40
$ws.Cells.Item(124, 11).Value2 = 12000
$ws.Cells.Item(124, 12).Value2 = 12000
$ws.Cells.Item(124, 13).Value2 = 12000
$ws.Cells.Item(124, 14).Value2 = '$/docena de paquetes'
$ws.Cells.Item(124, 15).Value2 = 'Provincia de Cautín'
$ws.Cells.Item(124, 16).Value2 = 1000
$ws.Cells.Item(124, 17).Value2 = 12
$ws.Cells.Item(124, 18).Value2 = 'Hortaliza'
